# Refresh the "cryptos" price/volume table with the latest scrape.
#
# Column D ("Price") and column E ("Volume(1h)") are stored as *text* in
# the workbook (e.g. "34.216.16", "  +0.72%  "), not numbers. Some of the
# new price strings (e.g. "226.69") are valid literal numbers, and a
# plain `Range.Value = "226.69"` assignment would make Excel coerce the
# cell to a numeric type, which would silently change the cell's stored
# type/format versus the source data. To avoid that, those values are
# entered with a leading apostrophe (forcing text, just like typing
# '226.69 into a cell), then the cell style is reset to "Normal" so no
# stray per-cell number-format is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '34.216.16'
$ws.Range('E2').Value = '  +0.72%  '
$ws.Range('D3').Value = '1.790.82'
$ws.Range('E3').Value = '  +0.44%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'226.69"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'31.97"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.92%  '
$ws.Range('E9').Value = '  +1.13%  '
$ws.Range('E10').Value = '  -1.97%  '
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = '2.049.96'
$ws.Range('D13').Value = "'11.14"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.61%  '
$ws.Range('D14').Value = '1.783.05'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').Value = '34.157.99'
$ws.Range('E15').Value = '  +0.48%  '
$ws.Range('D16').Value = "'0.621"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('E17').Value = '  +0.68%  '
$ws.Range('D18').Value = "'68.11"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.26%  '
$ws.Range('D19').Value = "'245.64"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('E20').Value = '  -0.64%  '
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').Value = "'10.83"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('D23').Value = "'4.12"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = "'2.05"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('D27').Value = "'16.37"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = "'1.23"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('D33').Value = "'3.64"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.85%  '
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').Value = '1.453.52'
$ws.Range('E35').Value = '  +4.27%  '
$ws.Range('E36').Value = '  -1.81%  '
$ws.Range('E37').Value = '  +7.75%  '
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('D40').Value = "'80.46"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.21%  '
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E42').Value = '  +0.87%  '
$ws.Range('E43').Value = '  +0.93%  '
$ws.Range('D44').Value = "'13.54"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('E45').Value = '  +2.43%  '
$ws.Range('E46').Value = '  +3.48%  '
$ws.Range('D47').Value = "'1.08"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').Value = '0.0₆0137'
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('D49').Value = '1.950.92'
$ws.Range('E49').Value = '  +0.66%  '
$ws.Range('D50').Value = "'105.94"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.13%  '
$ws.Range('E51').Value = '  +0.07%  '
